$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Run Date cell M1
$ws.Cells.Item(1, 13).Value = 45929

# Swap data rows 76 and 77
$ws.Cells.Item(76, 2).Value = 63793
$ws.Cells.Item(76, 5).Value = 267.11
$ws.Cells.Item(76, 6).Value = 1
$ws.Cells.Item(76, 7).Value = 251.25
$ws.Cells.Item(77, 2).Value = 44962
$ws.Cells.Item(77, 5).Value = 284.93
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0

# Swap data rows 82 and 83
$ws.Cells.Item(82, 2).Value = 44977
$ws.Cells.Item(82, 5).Value = 427.72
$ws.Cells.Item(82, 6).Value = 0
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(83, 2).Value = 63799
$ws.Cells.Item(83, 5).Value = 401
$ws.Cells.Item(83, 6).Value = 1
$ws.Cells.Item(83, 7).Value = 377.19

# Swap data rows 86 and 87
$ws.Cells.Item(86, 2).Value = 63794
$ws.Cells.Item(86, 5).Value = 668.78
$ws.Cells.Item(86, 6).Value = 7
$ws.Cells.Item(86, 7).Value = 4403.42
$ws.Cells.Item(87, 2).Value = 44964
$ws.Cells.Item(87, 5).Value = 713.34
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0

# Swap data rows 100 and 101
$ws.Cells.Item(100, 2).Value = 63838
$ws.Cells.Item(100, 5).Value = 150.09
$ws.Cells.Item(100, 6).Value = 5
$ws.Cells.Item(100, 7).Value = 705.95
$ws.Cells.Item(101, 2).Value = 46767
$ws.Cells.Item(101, 5).Value = 160.11
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0

# Swap data rows 189 and 190
$ws.Cells.Item(189, 2).Value = 48654
$ws.Cells.Item(189, 5).Value = 38.26
$ws.Cells.Item(189, 6).Value = 1
$ws.Cells.Item(189, 7).Value = 32.02
$ws.Cells.Item(190, 2).Value = 63902
$ws.Cells.Item(190, 5).Value = 34.04
$ws.Cells.Item(190, 6).Value = 2
$ws.Cells.Item(190, 7).Value = 64.04000000000001

# Swap data rows 234 and 235
$ws.Cells.Item(234, 2).Value = 64612
$ws.Cells.Item(234, 5).Value = 136.38
$ws.Cells.Item(234, 6).Value = 1
$ws.Cells.Item(234, 7).Value = 128.28
$ws.Cells.Item(235, 2).Value = 61866
$ws.Cells.Item(235, 5).Value = 145.48
$ws.Cells.Item(235, 6).Value = 0
$ws.Cells.Item(235, 7).Value = 0

# Swap data rows 366 and 367
$ws.Cells.Item(366, 2).Value = 62997
$ws.Cells.Item(366, 6).Value = 72
$ws.Cells.Item(366, 7).Value = 22020.48
$ws.Cells.Item(367, 2).Value = 57854
$ws.Cells.Item(367, 6).Value = 2
$ws.Cells.Item(367, 7).Value = 611.6799999999999

# Swap data rows 417 and 418
$ws.Cells.Item(417, 2).Value = 57882
$ws.Cells.Item(417, 5).Value = 58.2
$ws.Cells.Item(417, 6).Value = 0
$ws.Cells.Item(417, 7).Value = 0
$ws.Cells.Item(418, 2).Value = 63556
$ws.Cells.Item(418, 5).Value = 54.56
$ws.Cells.Item(418, 6).Value = 89
$ws.Cells.Item(418, 7).Value = 4567.48

# Swap data rows 485 and 486
$ws.Cells.Item(485, 2).Value = 57856
$ws.Cells.Item(485, 6).Value = 2
$ws.Cells.Item(485, 7).Value = 342.66
$ws.Cells.Item(486, 2).Value = 63007
$ws.Cells.Item(486, 6).Value = 984
$ws.Cells.Item(486, 7).Value = 168588.72

# Swap data rows 531 and 532
$ws.Cells.Item(531, 2).Value = 44198
$ws.Cells.Item(531, 5).Value = 38.52
$ws.Cells.Item(531, 6).Value = 0
$ws.Cells.Item(531, 7).Value = 0
$ws.Cells.Item(532, 2).Value = 63781
$ws.Cells.Item(532, 5).Value = 34.28
$ws.Cells.Item(532, 6).Value = 154
$ws.Cells.Item(532, 7).Value = 4964.96

# Swap data rows 535 and 536
$ws.Cells.Item(535, 2).Value = 54532
$ws.Cells.Item(535, 5).Value = 364.02
$ws.Cells.Item(535, 6).Value = 0
$ws.Cells.Item(535, 7).Value = 0
$ws.Cells.Item(536, 2).Value = 64190
$ws.Cells.Item(536, 5).Value = 341.28
$ws.Cells.Item(536, 6).Value = 16
$ws.Cells.Item(536, 7).Value = 5136.16

# Swap data rows 541 and 542
$ws.Cells.Item(541, 2).Value = 64191
$ws.Cells.Item(541, 5).Value = 341.28
$ws.Cells.Item(541, 6).Value = 2
$ws.Cells.Item(541, 7).Value = 642.02
$ws.Cells.Item(542, 2).Value = 54533
$ws.Cells.Item(542, 5).Value = 364.02
$ws.Cells.Item(542, 6).Value = 0
$ws.Cells.Item(542, 7).Value = 0

# Swap data rows 748 and 749
$ws.Cells.Item(748, 2).Value = 55658
$ws.Cells.Item(748, 5).Value = 801.25
$ws.Cells.Item(748, 6).Value = 0
$ws.Cells.Item(748, 7).Value = 0
$ws.Cells.Item(749, 2).Value = 64244
$ws.Cells.Item(749, 5).Value = 712.99
$ws.Cells.Item(749, 6).Value = 2
$ws.Cells.Item(749, 7).Value = 1341.3

# Swap data rows 776 and 777
$ws.Cells.Item(776, 2).Value = 46270
$ws.Cells.Item(776, 5).Value = 8.199999999999999
$ws.Cells.Item(776, 6).Value = 0
$ws.Cells.Item(776, 7).Value = 0
$ws.Cells.Item(777, 2).Value = 63810
$ws.Cells.Item(777, 5).Value = 7.28
$ws.Cells.Item(777, 6).Value = 64
$ws.Cells.Item(777, 7).Value = 438.4

# Swap data rows 778 and 779
$ws.Cells.Item(778, 2).Value = 46279
$ws.Cells.Item(778, 5).Value = 15.39
$ws.Cells.Item(778, 6).Value = 0
$ws.Cells.Item(778, 7).Value = 0
$ws.Cells.Item(779, 2).Value = 63816
$ws.Cells.Item(779, 5).Value = 13.69
$ws.Cells.Item(779, 6).Value = 228
$ws.Cells.Item(779, 7).Value = 2938.92

# Swap data rows 782 and 783
$ws.Cells.Item(782, 2).Value = 63812
$ws.Cells.Item(782, 5).Value = 7.12
$ws.Cells.Item(782, 6).Value = 2
$ws.Cells.Item(782, 7).Value = 13.4
$ws.Cells.Item(783, 2).Value = 46272
$ws.Cells.Item(783, 5).Value = 8
$ws.Cells.Item(783, 6).Value = 0
$ws.Cells.Item(783, 7).Value = 0

# Swap data rows 784 and 785
$ws.Cells.Item(784, 2).Value = 46266
$ws.Cells.Item(784, 5).Value = 19.76
$ws.Cells.Item(784, 6).Value = 0
$ws.Cells.Item(784, 7).Value = 0
$ws.Cells.Item(785, 2).Value = 63807
$ws.Cells.Item(785, 5).Value = 17.58
$ws.Cells.Item(785, 6).Value = 96
$ws.Cells.Item(785, 7).Value = 1586.88

# Swap data rows 788 and 789
$ws.Cells.Item(788, 2).Value = 63815
$ws.Cells.Item(788, 5).Value = 34.71
$ws.Cells.Item(788, 6).Value = 100
$ws.Cells.Item(788, 7).Value = 3264
$ws.Cells.Item(789, 2).Value = 46276
$ws.Cells.Item(789, 5).Value = 39.01
$ws.Cells.Item(789, 6).Value = 0
$ws.Cells.Item(789, 7).Value = 0

# Swap data rows 872 and 873
$ws.Cells.Item(872, 2).Value = 54098
$ws.Cells.Item(872, 5).Value = 134.46
$ws.Cells.Item(872, 6).Value = 0
$ws.Cells.Item(872, 7).Value = 0
$ws.Cells.Item(873, 2).Value = 64176
$ws.Cells.Item(873, 5).Value = 126.06
$ws.Cells.Item(873, 6).Value = 129
$ws.Cells.Item(873, 7).Value = 15296.82

# Swap data rows 884 and 885
$ws.Cells.Item(884, 2).Value = 64203
$ws.Cells.Item(884, 5).Value = 29.14
$ws.Cells.Item(884, 6).Value = 2
$ws.Cells.Item(884, 7).Value = 54.8
$ws.Cells.Item(885, 2).Value = 54894
$ws.Cells.Item(885, 5).Value = 32.74
$ws.Cells.Item(885, 6).Value = 0
$ws.Cells.Item(885, 7).Value = 0

# Swap data rows 887 and 888
$ws.Cells.Item(887, 2).Value = 64201
$ws.Cells.Item(887, 5).Value = 28.33
$ws.Cells.Item(887, 6).Value = 25
$ws.Cells.Item(887, 7).Value = 666
$ws.Cells.Item(888, 2).Value = 54892
$ws.Cells.Item(888, 5).Value = 31.83
$ws.Cells.Item(888, 6).Value = 0
$ws.Cells.Item(888, 7).Value = 0

# Swap data rows 896 and 897
$ws.Cells.Item(896, 2).Value = 64207
$ws.Cells.Item(896, 5).Value = 186.54
$ws.Cells.Item(896, 6).Value = 122
$ws.Cells.Item(896, 7).Value = 21407.34
$ws.Cells.Item(897, 2).Value = 54914
$ws.Cells.Item(897, 5).Value = 198.98
$ws.Cells.Item(897, 6).Value = 0
$ws.Cells.Item(897, 7).Value = 0
